# Insert a new price-record row for "Vega Modelo de Temuco - Albahaca" at
# row 169 (shifting the existing rows 169:271 down to 170:272) and fill the
# new row with the reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by one row, starting at row 169.
$ws.Rows.Item(169).Insert()

# Populate the newly-inserted row 169 with the new record.
$ws.Range("A169").Value = 10
$ws.Range("B169").Value = "Vega Modelo de Temuco"
$ws.Range("C169").Value = "La Araucanía"
$ws.Range("D169").Value = 44846
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112052
$ws.Range("G169").Value = "Albahaca"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 25
$ws.Range("K169").Value = 9000
$ws.Range("L169").Value = 9000
$ws.Range("M169").Value = 9000
$ws.Range("N169").Value = '$/paquete'
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 9000
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = "Hortaliza"
